$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to 2-decimal rounded precision (custom accuracy)
$ws.Range("B5").Value2 = 18.61
$ws.Range("C5").Value2 = 13.91
$ws.Range("D5").Value2 = 1.14
$ws.Range("E5").Value2 = 40.76
$ws.Range("F5").Value2 = 33.19
$ws.Range("G5").Value2 = 14.32
$ws.Range("H5").Value2 = 57.53
$ws.Range("I5").Value2 = 22.62
$ws.Range("J5").Value2 = 10.23
$ws.Range("K5").Value2 = 14.69
$ws.Range("L5").Value2 = 16.34
$ws.Range("M5").Value2 = 17.42
$ws.Range("N5").Value2 = 4.87
$ws.Range("O5").Value2 = 14.66
$ws.Range("P5").Value2 = 20.76
$ws.Range("Q5").Value2 = 12.45
$ws.Range("R5").Value2 = 0.49
$ws.Range("S5").Value2 = 0.71
$ws.Range("T5").Value2 = 215.86
$ws.Range("U5").Value2 = 41.01
$ws.Range("V5").Value2 = 13.53
$ws.Range("W5").Value2 = 27.48
$ws.Range("X5").Value2 = 14.47
$ws.Range("Y5").Value2 = 1.91
$ws.Range("Z5").Value2 = 28.44
$ws.Range("AA5").Value2 = 11.95
$ws.Range("AB5").Value2 = 10.62
$ws.Range("AC5").Value2 = 12.46
$ws.Range("AD5").Value2 = 17.22
$ws.Range("AE5").Value2 = 0.48
$ws.Range("AF5").Value2 = 52.33
$ws.Range("AG5").Value2 = 7.57
$ws.Range("AH5").Value2 = 16.92

# Column AB (28th column) width narrows from 8 to 7 characters
$ws.Columns.Item(28).ColumnWidth = 6.17

# Remove row 6 (data trimmed to 1000 rows -> fewer sample rows here)
$ws.Rows.Item(6).Delete()
